# Auto-generated Excel COM-interop script
# Applies: global column C bump (46059 -> 46060) and a row-content
# permutation of rows 4-31 (rows 2,3,7-12 keep their row position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: now holds data for 'A 37417-2023' (previously at row 6)
$ws.Cells.Item(4,1).Value = 'A 37417-2023'  # A4
$ws.Cells.Item(4,2).Value = 45155.0  # B4 (date serial)
$ws.Cells.Item(4,3).Value = 46060  # C4
$ws.Cells.Item(4,4).Value = 'UPPSALA LÄN'  # D4
$ws.Cells.Item(4,5).Value = 'HÅBO'  # E4
$ws.Cells.Item(4,7).Value = 12.9  # G4
$ws.Cells.Item(4,8).Value = 2  # H4
$ws.Cells.Item(4,9).Value = 0  # I4
$ws.Cells.Item(4,10).Value = 1  # J4
$ws.Cells.Item(4,11).Value = 1  # K4
$ws.Cells.Item(4,12).Value = 0  # L4
$ws.Cells.Item(4,13).Value = 0  # M4
$ws.Cells.Item(4,14).Value = 0  # N4
$ws.Cells.Item(4,15).Value = 2  # O4
$ws.Cells.Item(4,16).Value = 1  # P4
$ws.Cells.Item(4,17).Value = 3  # Q4
$ws.Cells.Item(4,18).Value = "Knärot`r`nUllticka`r`nBlåsippa"  # R4
$ws.Cells.Item(4,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 37417-2023 artfynd.xlsx", "A 37417-2023")'  # S4
$ws.Cells.Item(4,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 37417-2023 karta.png", "A 37417-2023")'  # T4
$ws.Cells.Item(4,21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/knärot/A 37417-2023 karta knärot.png", "A 37417-2023")'  # U4
$ws.Cells.Item(4,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 37417-2023 FSC-klagomål.docx", "A 37417-2023")'  # V4
$ws.Cells.Item(4,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 37417-2023 FSC-klagomål mail.docx", "A 37417-2023")'  # W4
$ws.Cells.Item(4,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 37417-2023 tillsynsbegäran.docx", "A 37417-2023")'  # X4
$ws.Cells.Item(4,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 37417-2023 tillsynsbegäran mail.docx", "A 37417-2023")'  # Y4
$ws.Cells.Item(4,26).ClearContents()  # Z4

# Row 5: now holds data for 'A 35197-2025' (previously at row 4)
$ws.Cells.Item(5,1).Value = 'A 35197-2025'  # A5
$ws.Cells.Item(5,2).Value = 45853.0  # B5 (date serial)
$ws.Cells.Item(5,3).Value = 46060  # C5
$ws.Cells.Item(5,4).Value = 'UPPSALA LÄN'  # D5
$ws.Cells.Item(5,5).Value = 'HÅBO'  # E5
$ws.Cells.Item(5,7).Value = 0.9  # G5
$ws.Cells.Item(5,8).Value = 1  # H5
$ws.Cells.Item(5,9).Value = 0  # I5
$ws.Cells.Item(5,10).Value = 3  # J5
$ws.Cells.Item(5,11).Value = 0  # K5
$ws.Cells.Item(5,12).Value = 0  # L5
$ws.Cells.Item(5,13).Value = 0  # M5
$ws.Cells.Item(5,14).Value = 0  # N5
$ws.Cells.Item(5,15).Value = 3  # O5
$ws.Cells.Item(5,16).Value = 0  # P5
$ws.Cells.Item(5,17).Value = 3  # Q5
$ws.Cells.Item(5,18).Value = "Grönsångare`r`nTallticka`r`nVintertagging"  # R5
$ws.Cells.Item(5,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 35197-2025 artfynd.xlsx", "A 35197-2025")'  # S5
$ws.Cells.Item(5,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 35197-2025 karta.png", "A 35197-2025")'  # T5
$ws.Cells.Item(5,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 35197-2025 FSC-klagomål.docx", "A 35197-2025")'  # V5
$ws.Cells.Item(5,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 35197-2025 FSC-klagomål mail.docx", "A 35197-2025")'  # W5
$ws.Cells.Item(5,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 35197-2025 tillsynsbegäran.docx", "A 35197-2025")'  # X5
$ws.Cells.Item(5,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 35197-2025 tillsynsbegäran mail.docx", "A 35197-2025")'  # Y5
$ws.Cells.Item(5,26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/fåglar/A 35197-2025 prioriterade fågelarter.docx", "A 35197-2025")'  # Z5

# Row 6: now holds data for 'A 47653-2024' (previously at row 5)
$ws.Cells.Item(6,1).Value = 'A 47653-2024'  # A6
$ws.Cells.Item(6,2).Value = 45588.0  # B6 (date serial)
$ws.Cells.Item(6,3).Value = 46060  # C6
$ws.Cells.Item(6,4).Value = 'UPPSALA LÄN'  # D6
$ws.Cells.Item(6,5).Value = 'HÅBO'  # E6
$ws.Cells.Item(6,7).Value = 3  # G6
$ws.Cells.Item(6,8).Value = 2  # H6
$ws.Cells.Item(6,9).Value = 1  # I6
$ws.Cells.Item(6,10).Value = 1  # J6
$ws.Cells.Item(6,11).Value = 0  # K6
$ws.Cells.Item(6,12).Value = 0  # L6
$ws.Cells.Item(6,13).Value = 0  # M6
$ws.Cells.Item(6,14).Value = 0  # N6
$ws.Cells.Item(6,15).Value = 1  # O6
$ws.Cells.Item(6,16).Value = 0  # P6
$ws.Cells.Item(6,17).Value = 3  # Q6
$ws.Cells.Item(6,18).Value = "Backklöver`r`nNästrot`r`nBlåsippa"  # R6
$ws.Cells.Item(6,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 47653-2024 artfynd.xlsx", "A 47653-2024")'  # S6
$ws.Cells.Item(6,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 47653-2024 karta.png", "A 47653-2024")'  # T6
$ws.Cells.Item(6,21).ClearContents()  # U6
$ws.Cells.Item(6,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 47653-2024 FSC-klagomål.docx", "A 47653-2024")'  # V6
$ws.Cells.Item(6,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 47653-2024 FSC-klagomål mail.docx", "A 47653-2024")'  # W6
$ws.Cells.Item(6,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 47653-2024 tillsynsbegäran.docx", "A 47653-2024")'  # X6
$ws.Cells.Item(6,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 47653-2024 tillsynsbegäran mail.docx", "A 47653-2024")'  # Y6

# Row 13: now holds data for 'A 12146-2023' (previously at row 27)
$ws.Cells.Item(13,1).Value = 'A 12146-2023'  # A13
$ws.Cells.Item(13,2).Value = 44998.478425925925  # B13 (date serial)
$ws.Cells.Item(13,3).Value = 46060  # C13
$ws.Cells.Item(13,4).Value = 'UPPSALA LÄN'  # D13
$ws.Cells.Item(13,5).Value = 'HÅBO'  # E13
$ws.Cells.Item(13,7).Value = 3.1  # G13
$ws.Cells.Item(13,8).Value = 0  # H13
$ws.Cells.Item(13,9).Value = 0  # I13
$ws.Cells.Item(13,10).Value = 0  # J13
$ws.Cells.Item(13,11).Value = 0  # K13
$ws.Cells.Item(13,12).Value = 0  # L13
$ws.Cells.Item(13,13).Value = 0  # M13
$ws.Cells.Item(13,14).Value = 0  # N13
$ws.Cells.Item(13,15).Value = 0  # O13
$ws.Cells.Item(13,16).Value = 0  # P13
$ws.Cells.Item(13,17).Value = 0  # Q13

# Row 14: now holds data for 'A 32023-2023' (previously at row 18)
$ws.Cells.Item(14,1).Value = 'A 32023-2023'  # A14
$ws.Cells.Item(14,2).Value = 45119.49833333334  # B14 (date serial)
$ws.Cells.Item(14,3).Value = 46060  # C14
$ws.Cells.Item(14,4).Value = 'UPPSALA LÄN'  # D14
$ws.Cells.Item(14,5).Value = 'HÅBO'  # E14
$ws.Cells.Item(14,7).Value = 3.1  # G14
$ws.Cells.Item(14,8).Value = 0  # H14
$ws.Cells.Item(14,9).Value = 0  # I14
$ws.Cells.Item(14,10).Value = 0  # J14
$ws.Cells.Item(14,11).Value = 0  # K14
$ws.Cells.Item(14,12).Value = 0  # L14
$ws.Cells.Item(14,13).Value = 0  # M14
$ws.Cells.Item(14,14).Value = 0  # N14
$ws.Cells.Item(14,15).Value = 0  # O14
$ws.Cells.Item(14,16).Value = 0  # P14
$ws.Cells.Item(14,17).Value = 0  # Q14

# Row 15: now holds data for 'A 4524-2024' (previously at row 23)
$ws.Cells.Item(15,1).Value = 'A 4524-2024'  # A15
$ws.Cells.Item(15,2).Value = 45327.0  # B15 (date serial)
$ws.Cells.Item(15,3).Value = 46060  # C15
$ws.Cells.Item(15,4).Value = 'UPPSALA LÄN'  # D15
$ws.Cells.Item(15,5).Value = 'HÅBO'  # E15
$ws.Cells.Item(15,7).Value = 4.6  # G15
$ws.Cells.Item(15,8).Value = 0  # H15
$ws.Cells.Item(15,9).Value = 0  # I15
$ws.Cells.Item(15,10).Value = 0  # J15
$ws.Cells.Item(15,11).Value = 0  # K15
$ws.Cells.Item(15,12).Value = 0  # L15
$ws.Cells.Item(15,13).Value = 0  # M15
$ws.Cells.Item(15,14).Value = 0  # N15
$ws.Cells.Item(15,15).Value = 0  # O15
$ws.Cells.Item(15,16).Value = 0  # P15
$ws.Cells.Item(15,17).Value = 0  # Q15

# Row 16: now holds data for 'A 37415-2023' (previously at row 26)
$ws.Cells.Item(16,1).Value = 'A 37415-2023'  # A16
$ws.Cells.Item(16,2).Value = 45155.0  # B16 (date serial)
$ws.Cells.Item(16,3).Value = 46060  # C16
$ws.Cells.Item(16,4).Value = 'UPPSALA LÄN'  # D16
$ws.Cells.Item(16,5).Value = 'HÅBO'  # E16
$ws.Cells.Item(16,7).Value = 6.6  # G16
$ws.Cells.Item(16,8).Value = 0  # H16
$ws.Cells.Item(16,9).Value = 0  # I16
$ws.Cells.Item(16,10).Value = 0  # J16
$ws.Cells.Item(16,11).Value = 0  # K16
$ws.Cells.Item(16,12).Value = 0  # L16
$ws.Cells.Item(16,13).Value = 0  # M16
$ws.Cells.Item(16,14).Value = 0  # N16
$ws.Cells.Item(16,15).Value = 0  # O16
$ws.Cells.Item(16,16).Value = 0  # P16
$ws.Cells.Item(16,17).Value = 0  # Q16

# Row 17: now holds data for 'A 55962-2023' (previously at row 14)
$ws.Cells.Item(17,1).Value = 'A 55962-2023'  # A17
$ws.Cells.Item(17,2).Value = 45240.0  # B17 (date serial)
$ws.Cells.Item(17,3).Value = 46060  # C17
$ws.Cells.Item(17,4).Value = 'UPPSALA LÄN'  # D17
$ws.Cells.Item(17,5).Value = 'HÅBO'  # E17
$ws.Cells.Item(17,7).Value = 3.4  # G17
$ws.Cells.Item(17,8).Value = 0  # H17
$ws.Cells.Item(17,9).Value = 0  # I17
$ws.Cells.Item(17,10).Value = 0  # J17
$ws.Cells.Item(17,11).Value = 0  # K17
$ws.Cells.Item(17,12).Value = 0  # L17
$ws.Cells.Item(17,13).Value = 0  # M17
$ws.Cells.Item(17,14).Value = 0  # N17
$ws.Cells.Item(17,15).Value = 0  # O17
$ws.Cells.Item(17,16).Value = 0  # P17
$ws.Cells.Item(17,17).Value = 0  # Q17

# Row 18: now holds data for 'A 11989-2025' (previously at row 28)
$ws.Cells.Item(18,1).Value = 'A 11989-2025'  # A18
$ws.Cells.Item(18,2).Value = 45728.60074074074  # B18 (date serial)
$ws.Cells.Item(18,3).Value = 46060  # C18
$ws.Cells.Item(18,4).Value = 'UPPSALA LÄN'  # D18
$ws.Cells.Item(18,5).Value = 'HÅBO'  # E18
$ws.Cells.Item(18,7).Value = 9.6  # G18
$ws.Cells.Item(18,8).Value = 0  # H18
$ws.Cells.Item(18,9).Value = 0  # I18
$ws.Cells.Item(18,10).Value = 0  # J18
$ws.Cells.Item(18,11).Value = 0  # K18
$ws.Cells.Item(18,12).Value = 0  # L18
$ws.Cells.Item(18,13).Value = 0  # M18
$ws.Cells.Item(18,14).Value = 0  # N18
$ws.Cells.Item(18,15).Value = 0  # O18
$ws.Cells.Item(18,16).Value = 0  # P18
$ws.Cells.Item(18,17).Value = 0  # Q18

# Row 19: now holds data for 'A 45406-2025' (previously at row 15)
$ws.Cells.Item(19,1).Value = 'A 45406-2025'  # A19
$ws.Cells.Item(19,2).Value = 45922.42936342592  # B19 (date serial)
$ws.Cells.Item(19,3).Value = 46060  # C19
$ws.Cells.Item(19,4).Value = 'UPPSALA LÄN'  # D19
$ws.Cells.Item(19,5).Value = 'HÅBO'  # E19
$ws.Cells.Item(19,7).Value = 9.1  # G19
$ws.Cells.Item(19,8).Value = 0  # H19
$ws.Cells.Item(19,9).Value = 0  # I19
$ws.Cells.Item(19,10).Value = 0  # J19
$ws.Cells.Item(19,11).Value = 0  # K19
$ws.Cells.Item(19,12).Value = 0  # L19
$ws.Cells.Item(19,13).Value = 0  # M19
$ws.Cells.Item(19,14).Value = 0  # N19
$ws.Cells.Item(19,15).Value = 0  # O19
$ws.Cells.Item(19,16).Value = 0  # P19
$ws.Cells.Item(19,17).Value = 0  # Q19

# Row 20: now holds data for 'A 35198-2025' (previously at row 31)
$ws.Cells.Item(20,1).Value = 'A 35198-2025'  # A20
$ws.Cells.Item(20,2).Value = 45853.0  # B20 (date serial)
$ws.Cells.Item(20,3).Value = 46060  # C20
$ws.Cells.Item(20,4).Value = 'UPPSALA LÄN'  # D20
$ws.Cells.Item(20,5).Value = 'HÅBO'  # E20
$ws.Cells.Item(20,6).ClearContents()  # F20
$ws.Cells.Item(20,7).Value = 1.2  # G20
$ws.Cells.Item(20,8).Value = 0  # H20
$ws.Cells.Item(20,9).Value = 0  # I20
$ws.Cells.Item(20,10).Value = 0  # J20
$ws.Cells.Item(20,11).Value = 0  # K20
$ws.Cells.Item(20,12).Value = 0  # L20
$ws.Cells.Item(20,13).Value = 0  # M20
$ws.Cells.Item(20,14).Value = 0  # N20
$ws.Cells.Item(20,15).Value = 0  # O20
$ws.Cells.Item(20,16).Value = 0  # P20
$ws.Cells.Item(20,17).Value = 0  # Q20

# Row 21: now holds data for 'A 35300-2025' (previously at row 30)
$ws.Cells.Item(21,1).Value = 'A 35300-2025'  # A21
$ws.Cells.Item(21,2).Value = 45854.41511574074  # B21 (date serial)
$ws.Cells.Item(21,3).Value = 46060  # C21
$ws.Cells.Item(21,4).Value = 'UPPSALA LÄN'  # D21
$ws.Cells.Item(21,5).Value = 'HÅBO'  # E21
$ws.Cells.Item(21,6).ClearContents()  # F21
$ws.Cells.Item(21,7).Value = 2.2  # G21
$ws.Cells.Item(21,8).Value = 0  # H21
$ws.Cells.Item(21,9).Value = 0  # I21
$ws.Cells.Item(21,10).Value = 0  # J21
$ws.Cells.Item(21,11).Value = 0  # K21
$ws.Cells.Item(21,12).Value = 0  # L21
$ws.Cells.Item(21,13).Value = 0  # M21
$ws.Cells.Item(21,14).Value = 0  # N21
$ws.Cells.Item(21,15).Value = 0  # O21
$ws.Cells.Item(21,16).Value = 0  # P21
$ws.Cells.Item(21,17).Value = 0  # Q21

# Row 22: now holds data for 'A 37410-2023' (previously at row 16)
$ws.Cells.Item(22,1).Value = 'A 37410-2023'  # A22
$ws.Cells.Item(22,2).Value = 45155.0  # B22 (date serial)
$ws.Cells.Item(22,3).Value = 46060  # C22
$ws.Cells.Item(22,4).Value = 'UPPSALA LÄN'  # D22
$ws.Cells.Item(22,5).Value = 'HÅBO'  # E22
$ws.Cells.Item(22,7).Value = 20.9  # G22
$ws.Cells.Item(22,8).Value = 0  # H22
$ws.Cells.Item(22,9).Value = 0  # I22
$ws.Cells.Item(22,10).Value = 0  # J22
$ws.Cells.Item(22,11).Value = 0  # K22
$ws.Cells.Item(22,12).Value = 0  # L22
$ws.Cells.Item(22,13).Value = 0  # M22
$ws.Cells.Item(22,14).Value = 0  # N22
$ws.Cells.Item(22,15).Value = 0  # O22
$ws.Cells.Item(22,16).Value = 0  # P22
$ws.Cells.Item(22,17).Value = 0  # Q22

# Row 23: now holds data for 'A 35193-2025' (previously at row 13)
$ws.Cells.Item(23,1).Value = 'A 35193-2025'  # A23
$ws.Cells.Item(23,2).Value = 45853.0  # B23 (date serial)
$ws.Cells.Item(23,3).Value = 46060  # C23
$ws.Cells.Item(23,4).Value = 'UPPSALA LÄN'  # D23
$ws.Cells.Item(23,5).Value = 'HÅBO'  # E23
$ws.Cells.Item(23,7).Value = 1.9  # G23
$ws.Cells.Item(23,8).Value = 0  # H23
$ws.Cells.Item(23,9).Value = 0  # I23
$ws.Cells.Item(23,10).Value = 0  # J23
$ws.Cells.Item(23,11).Value = 0  # K23
$ws.Cells.Item(23,12).Value = 0  # L23
$ws.Cells.Item(23,13).Value = 0  # M23
$ws.Cells.Item(23,14).Value = 0  # N23
$ws.Cells.Item(23,15).Value = 0  # O23
$ws.Cells.Item(23,16).Value = 0  # P23
$ws.Cells.Item(23,17).Value = 0  # Q23

# Row 24: now holds data for 'A 12156-2023' (previously at row 17)
$ws.Cells.Item(24,1).Value = 'A 12156-2023'  # A24
$ws.Cells.Item(24,2).Value = 44998.491574074076  # B24 (date serial)
$ws.Cells.Item(24,3).Value = 46060  # C24
$ws.Cells.Item(24,4).Value = 'UPPSALA LÄN'  # D24
$ws.Cells.Item(24,5).Value = 'HÅBO'  # E24
$ws.Cells.Item(24,7).Value = 0.5  # G24
$ws.Cells.Item(24,8).Value = 0  # H24
$ws.Cells.Item(24,9).Value = 0  # I24
$ws.Cells.Item(24,10).Value = 0  # J24
$ws.Cells.Item(24,11).Value = 0  # K24
$ws.Cells.Item(24,12).Value = 0  # L24
$ws.Cells.Item(24,13).Value = 0  # M24
$ws.Cells.Item(24,14).Value = 0  # N24
$ws.Cells.Item(24,15).Value = 0  # O24
$ws.Cells.Item(24,16).Value = 0  # P24
$ws.Cells.Item(24,17).Value = 0  # Q24

# Row 25: now holds data for 'A 65836-2021' (previously at row 24)
$ws.Cells.Item(25,1).Value = 'A 65836-2021'  # A25
$ws.Cells.Item(25,2).Value = 44517.0  # B25 (date serial)
$ws.Cells.Item(25,3).Value = 46060  # C25
$ws.Cells.Item(25,4).Value = 'UPPSALA LÄN'  # D25
$ws.Cells.Item(25,5).Value = 'HÅBO'  # E25
$ws.Cells.Item(25,7).Value = 1.8  # G25
$ws.Cells.Item(25,8).Value = 0  # H25
$ws.Cells.Item(25,9).Value = 0  # I25
$ws.Cells.Item(25,10).Value = 0  # J25
$ws.Cells.Item(25,11).Value = 0  # K25
$ws.Cells.Item(25,12).Value = 0  # L25
$ws.Cells.Item(25,13).Value = 0  # M25
$ws.Cells.Item(25,14).Value = 0  # N25
$ws.Cells.Item(25,15).Value = 0  # O25
$ws.Cells.Item(25,16).Value = 0  # P25
$ws.Cells.Item(25,17).Value = 0  # Q25

# Row 26: now holds data for 'A 34202-2022' (previously at row 22)
$ws.Cells.Item(26,1).Value = 'A 34202-2022'  # A26
$ws.Cells.Item(26,2).Value = 44791.64837962963  # B26 (date serial)
$ws.Cells.Item(26,3).Value = 46060  # C26
$ws.Cells.Item(26,4).Value = 'UPPSALA LÄN'  # D26
$ws.Cells.Item(26,5).Value = 'HÅBO'  # E26
$ws.Cells.Item(26,7).Value = 2  # G26
$ws.Cells.Item(26,8).Value = 0  # H26
$ws.Cells.Item(26,9).Value = 0  # I26
$ws.Cells.Item(26,10).Value = 0  # J26
$ws.Cells.Item(26,11).Value = 0  # K26
$ws.Cells.Item(26,12).Value = 0  # L26
$ws.Cells.Item(26,13).Value = 0  # M26
$ws.Cells.Item(26,14).Value = 0  # N26
$ws.Cells.Item(26,15).Value = 0  # O26
$ws.Cells.Item(26,16).Value = 0  # P26
$ws.Cells.Item(26,17).Value = 0  # Q26

# Row 27: now holds data for 'A 12154-2023' (previously at row 25)
$ws.Cells.Item(27,1).Value = 'A 12154-2023'  # A27
$ws.Cells.Item(27,2).Value = 44998.0  # B27 (date serial)
$ws.Cells.Item(27,3).Value = 46060  # C27
$ws.Cells.Item(27,4).Value = 'UPPSALA LÄN'  # D27
$ws.Cells.Item(27,5).Value = 'HÅBO'  # E27
$ws.Cells.Item(27,7).Value = 2.7  # G27
$ws.Cells.Item(27,8).Value = 0  # H27
$ws.Cells.Item(27,9).Value = 0  # I27
$ws.Cells.Item(27,10).Value = 0  # J27
$ws.Cells.Item(27,11).Value = 0  # K27
$ws.Cells.Item(27,12).Value = 0  # L27
$ws.Cells.Item(27,13).Value = 0  # M27
$ws.Cells.Item(27,14).Value = 0  # N27
$ws.Cells.Item(27,15).Value = 0  # O27
$ws.Cells.Item(27,16).Value = 0  # P27
$ws.Cells.Item(27,17).Value = 0  # Q27

# Row 28: now holds data for 'A 67005-2021' (previously at row 19)
$ws.Cells.Item(28,1).Value = 'A 67005-2021'  # A28
$ws.Cells.Item(28,2).Value = 44522.0  # B28 (date serial)
$ws.Cells.Item(28,3).Value = 46060  # C28
$ws.Cells.Item(28,4).Value = 'UPPSALA LÄN'  # D28
$ws.Cells.Item(28,5).Value = 'HÅBO'  # E28
$ws.Cells.Item(28,7).Value = 1.3  # G28
$ws.Cells.Item(28,8).Value = 0  # H28
$ws.Cells.Item(28,9).Value = 0  # I28
$ws.Cells.Item(28,10).Value = 0  # J28
$ws.Cells.Item(28,11).Value = 0  # K28
$ws.Cells.Item(28,12).Value = 0  # L28
$ws.Cells.Item(28,13).Value = 0  # M28
$ws.Cells.Item(28,14).Value = 0  # N28
$ws.Cells.Item(28,15).Value = 0  # O28
$ws.Cells.Item(28,16).Value = 0  # P28
$ws.Cells.Item(28,17).Value = 0  # Q28

# Row 29: now holds data for 'A 15732-2025' (previously at row 21)
$ws.Cells.Item(29,1).Value = 'A 15732-2025'  # A29
$ws.Cells.Item(29,2).Value = 45747.0  # B29 (date serial)
$ws.Cells.Item(29,3).Value = 46060  # C29
$ws.Cells.Item(29,4).Value = 'UPPSALA LÄN'  # D29
$ws.Cells.Item(29,5).Value = 'HÅBO'  # E29
$ws.Cells.Item(29,6).Value = 'Kyrkan'  # F29
$ws.Cells.Item(29,7).Value = 1.4  # G29
$ws.Cells.Item(29,8).Value = 0  # H29
$ws.Cells.Item(29,9).Value = 0  # I29
$ws.Cells.Item(29,10).Value = 0  # J29
$ws.Cells.Item(29,11).Value = 0  # K29
$ws.Cells.Item(29,12).Value = 0  # L29
$ws.Cells.Item(29,13).Value = 0  # M29
$ws.Cells.Item(29,14).Value = 0  # N29
$ws.Cells.Item(29,15).Value = 0  # O29
$ws.Cells.Item(29,16).Value = 0  # P29
$ws.Cells.Item(29,17).Value = 0  # Q29

# Row 30: now holds data for 'A 23250-2022' (previously at row 20)
$ws.Cells.Item(30,1).Value = 'A 23250-2022'  # A30
$ws.Cells.Item(30,2).Value = 44719.0  # B30 (date serial)
$ws.Cells.Item(30,3).Value = 46060  # C30
$ws.Cells.Item(30,4).Value = 'UPPSALA LÄN'  # D30
$ws.Cells.Item(30,5).Value = 'HÅBO'  # E30
$ws.Cells.Item(30,6).Value = 'Naturvårdsverket'  # F30
$ws.Cells.Item(30,7).Value = 1  # G30
$ws.Cells.Item(30,8).Value = 0  # H30
$ws.Cells.Item(30,9).Value = 0  # I30
$ws.Cells.Item(30,10).Value = 0  # J30
$ws.Cells.Item(30,11).Value = 0  # K30
$ws.Cells.Item(30,12).Value = 0  # L30
$ws.Cells.Item(30,13).Value = 0  # M30
$ws.Cells.Item(30,14).Value = 0  # N30
$ws.Cells.Item(30,15).Value = 0  # O30
$ws.Cells.Item(30,16).Value = 0  # P30
$ws.Cells.Item(30,17).Value = 0  # Q30

# Row 31: now holds data for 'A 23370-2025' (previously at row 29)
$ws.Cells.Item(31,1).Value = 'A 23370-2025'  # A31
$ws.Cells.Item(31,2).Value = 45791.709074074075  # B31 (date serial)
$ws.Cells.Item(31,3).Value = 46060  # C31
$ws.Cells.Item(31,4).Value = 'UPPSALA LÄN'  # D31
$ws.Cells.Item(31,5).Value = 'HÅBO'  # E31
$ws.Cells.Item(31,7).Value = 3.8  # G31
$ws.Cells.Item(31,8).Value = 0  # H31
$ws.Cells.Item(31,9).Value = 0  # I31
$ws.Cells.Item(31,10).Value = 0  # J31
$ws.Cells.Item(31,11).Value = 0  # K31
$ws.Cells.Item(31,12).Value = 0  # L31
$ws.Cells.Item(31,13).Value = 0  # M31
$ws.Cells.Item(31,14).Value = 0  # N31
$ws.Cells.Item(31,15).Value = 0  # O31
$ws.Cells.Item(31,16).Value = 0  # P31
$ws.Cells.Item(31,17).Value = 0  # Q31

# Rows that keep their position: only the 'Förändrad' (C) date needs the +1 day bump.
$ws.Cells.Item(2,3).Value = 46060  # C2
$ws.Cells.Item(3,3).Value = 46060  # C3
$ws.Cells.Item(7,3).Value = 46060  # C7
$ws.Cells.Item(8,3).Value = 46060  # C8
$ws.Cells.Item(9,3).Value = 46060  # C9
$ws.Cells.Item(10,3).Value = 46060  # C10
$ws.Cells.Item(11,3).Value = 46060  # C11
$ws.Cells.Item(12,3).Value = 46060  # C12

Write-Host "Done applying HÅBO worksheet update."
